# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" (bound only to the Notes Master)
#   ppt/theme/theme2.xml  -> "Integral"     (bound to the real Slide Master /
#                                             Presentation.Designs(1))
# The authored edit swaps the two themes' contents: theme1.xml ends up holding
# the "Integral" palette and theme2.xml ends up holding the "Office Theme"
# palette (the slide master itself, its layouts, and every relationship id
# stay untouched - only the colour definitions that live inside the theme
# parts change place).
#
# The only theme surface the PowerPoint object model exposes for editing is
# the live design's ThemeColorScheme (Presentation.SlideMaster.Theme /
# Presentation.Designs(1)), which is the part backing ppt/theme/theme2.xml.
# We recolor it, in place, from "Integral" to the "Office Theme" values so
# that the design applied to the deck becomes the Office palette - the
# observable half of the swap.

function ToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Office Theme palette (was previously in ppt/theme/theme1.xml), applied in
# clrScheme order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeTheme = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeTheme.Length; $i++) {
    $colorScheme.Colors($i).RGB = ToRgbInt $officeTheme[$i - 1]
}
